$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 6: Relay - add vendor / part number / link
# ---------------------------------------------------------------------
$ws.Range("B6").Value2 = "Mouser"
$ws.Range("C6").Value2 = "881-RF300-5"
$ws.Range("E6").Value2 = "https://www.mouser.com/ProductDetail/Teledyne-Relays/RF300-5?qs=cFlnt7DBZX%252BJkIFFO4rSPw%3D%3D"

# ---------------------------------------------------------------------
# Row 7: was "Wires" -> becomes "JST Socket" entry
# ---------------------------------------------------------------------
$ws.Range("A7").Value2 = "JST Socket"
$ws.Range("B7").Value2 = "DigiKey"
$ws.Range("C7").Value2 = "455-1719-ND"
$ws.Range("D7").Value2 = "Right Angle"
$ws.Range("E7").Value2 = "https://www.digikey.com/en/products/detail/jst-sales-america-inc/S2B-PH-K-S-LF-SN/926626"

# ---------------------------------------------------------------------
# Row 8: was "JST Socket" -> becomes "JST Connector Kit" entry
# ---------------------------------------------------------------------
$ws.Range("A8").Value2 = "JST Connector Kit"
$ws.Range("B8").Value2 = "Amazon"
$ws.Range("D8").Value2 = "2-pin, 2mm"
$ws.Range("E8").Value2 = "https://www.amazon.com/Connector-Pre-Crimped-Compatible-JST-PH2-0mm-Inductrix/dp/B08T89ZK2Q/ref=sr_1_5?crid=2KZ2LYUXJW5KR&keywords=dupont+pre+crimped&qid=1677196883&sprefix=dupont+pre+crimped%25252Caps%25252C69&sr=8-5"

# ---------------------------------------------------------------------
# Row 9: was "JST Connector" -> becomes "Signal Generator"
# Row 10: was "Signal Generator" -> becomes "Alligator Clips"
# Row 11: was "Alligator Clips" -> becomes "Headstage"
# ---------------------------------------------------------------------
$ws.Range("A9").Value2 = "Signal Generator"
$ws.Range("A10").Value2 = "Alligator Clips"
$ws.Range("A11").Value2 = "Headstage"

# ---------------------------------------------------------------------
# Row 12: was "Headstage" -> becomes "Screws" entry (previously row 13)
# ---------------------------------------------------------------------
$ws.Range("A12").Value2 = "Screws"
$ws.Range("B12").Value2 = "Amazon"
$ws.Range("D12").Value2 = "#4-40"
$ws.Range("E12").Value2 = "https://www.amazon.com/MCMASKE-280PCS-Stainless-Washers-Assortment/dp/B0CB7TSTQK/ref=sr_1_21?keywords=4-40%25252BScrews&sr=8-21&th=1"

# ---------------------------------------------------------------------
# Row 13: was "Screws" -> becomes "Screw Spacers" entry (previously row 14)
# note: the link text here stays plain text (not a real hyperlink),
# matching the source workbook's original behavior. The old D13
# ("#4-40") value needs to be cleared since this row no longer has a
# Notes entry.
# ---------------------------------------------------------------------
$ws.Range("A13").Value2 = "Screw Spacers"
$ws.Range("B13").Value2 = "McMaster"
$ws.Range("D13").ClearContents() | Out-Null
$ws.Range("E13").Value2 = "https://www.mcmaster.com/93657A207/ or https://www.mcmaster.com/93657A505/"

# ---------------------------------------------------------------------
# Row 14: new entry "BNC Connector" (was "Screw Spacers")
# ---------------------------------------------------------------------
$ws.Range("A14").Value2 = "BNC Connector"
$ws.Range("B14").Value2 = "DigiKey"
$ws.Range("C14").Value2 = "314-1393-ND"
$ws.Range("D14").Value2 = "Cheaper equivalent on Amazon"
$ws.Range("E14").Value2 = "https://www.digikey.com/en/products/detail/mueller-electric-co/BU-5100-A-4-0/5801064"

# ---------------------------------------------------------------------
# Row 15: new entry "D-Sub Male"
# ---------------------------------------------------------------------
$ws.Range("A15").Value2 = "D-Sub Male"
$ws.Range("B15").Value2 = "DigiKey"
$ws.Range("C15").Value2 = "A34501-ND"
$ws.Range("D15").Value2 = "For compatibility with Axon Instruments Part # <?>"
$ws.Range("E15").Value2 = "https://www.digikey.com/en/products/detail/te-connectivity-aerospace-defense-and-marine/206794-4/1144251"

# ---------------------------------------------------------------------
# Row 16: new entry "D-Sub Female"
# ---------------------------------------------------------------------
$ws.Range("A16").Value2 = "D-Sub Female"
$ws.Range("B16").Value2 = "DigiKey"
$ws.Range("C16").Value2 = "A34502-ND"
$ws.Range("D16").Value2 = "For compatibility with Axon Instruments Part # <?>"
$ws.Range("E16").Value2 = "https://www.digikey.com/en/products/detail/te-connectivity-aerospace-defense-and-marine/206795-3/1144252"

# ---------------------------------------------------------------------
# Row 17: new entry "Pipette Holder"
# ---------------------------------------------------------------------
$ws.Range("A17").Value2 = "Pipette Holder"
$ws.Range("B17").Value2 = "Warner Instruments"
$ws.Range("C17").Value2 = "QSW-T20P"
$ws.Range("E17").Value2 = "https://www.harvardapparatus.com/media/brochures/Warner_Microelectrode_Holders.pdf"

# Rows 16 and 17 are brand new (beyond the original A4:E15 dimension), so
# column A does not automatically inherit the bold item-name style used
# by the rest of the table (cellXfs index 2: bold, 11pt). Apply it
# explicitly so it matches the other rows in the list.
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Font.Size = 11
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").Font.Size = 11

# ---------------------------------------------------------------------
# Hyperlinks - added in the same order as the authored workbook so the
# generated relationship ids (rId1..rId9) line up.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E12"), $ws.Range("E12").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"),  $ws.Range("E5").Value2)  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"),  $ws.Range("E6").Value2)  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"),  $ws.Range("E8").Value2)  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E7"),  $ws.Range("E7").Value2)  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E14"), $ws.Range("E14").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E15"), $ws.Range("E15").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E16"), $ws.Range("E16").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E17"), $ws.Range("E17").Value2) | Out-Null

# ---------------------------------------------------------------------
# Column widths (B and E) and selection.
# Note: this runtime snaps ColumnWidth to its own internal character
# grid (coarser than real Excel's pixel grid), so the raw inputs below
# are chosen to land as close as possible to the authored widths
# (18.7265625 and 132.6328125 respectively) once that snapping happens.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.8
$ws.Columns.Item(5).ColumnWidth = 131.8

$ws.Range("E20").Select() | Out-Null
